$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 17:56"

# Swap Bahamas / Estonia rows (row 136 becomes Bahamas, row 137 becomes Estonia)
# and refresh their data, plus update all other changed country stats.

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6491391
$ws.Range("C4").Value = 5816
$ws.Range("E4").Value = 2538600
$ws.Range("G4").Value = 123
$ws.Range("H4").Value = 193657

# Row 22 - Italia
$ws.Range("B22").Value = 280153
$ws.Range("C22").Value = 1370
$ws.Range("D22").Value = 210801
$ws.Range("E22").Value = 33789
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 35563

# Row 24 - Alemania
$ws.Range("B24").Value = 254615
$ws.Range("C24").Value = 990
$ws.Range("E24").Value = 18208

# Row 29 - Canada
$ws.Range("B29").Value = 132680
$ws.Range("C29").Value = 538
$ws.Range("D29").Value = 116900
$ws.Range("E29").Value = 6634

# Row 31 - Catar
$ws.Range("B31").Value = 120579
$ws.Range("C31").Value = 231
$ws.Range("D31").Value = 117497
$ws.Range("E31").Value = 2877

# Row 45 - Emiratos Arabes Unidos
$ws.Range("B45").Value = 75098
$ws.Range("C45").Value = 644
$ws.Range("D45").Value = 66943
$ws.Range("E45").Value = 7764
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 391

# Row 47 - Bielorrusia
$ws.Range("D47").Value = 71916
$ws.Range("E47").Value = 571

# Row 53 - Singapur
$ws.Range("D53").Value = 56461
$ws.Range("E53").Value = 603

# Row 65 - Moldavia
$ws.Range("B65").Value = 40556
$ws.Range("C65").Value = 501
$ws.Range("E65").Value = 10891
$ws.Range("G65").Value = 13
$ws.Range("H65").Value = 1087

# Row 91 - Grecia
$ws.Range("B91").Value = 11832
$ws.Range("C91").Value = 169
$ws.Range("E91").Value = 7738
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 290

# Row 93 - Albania
$ws.Range("B93").Value = 10553
$ws.Range("C93").Value = 147
$ws.Range("D93").Value = 6239
$ws.Range("E93").Value = 3995

# Row 133 - Sri Lanka
$ws.Range("B133").Value = 3131
$ws.Range("C133").Value = 8
$ws.Range("D133").Value = 2935
$ws.Range("E133").Value = 184

# Row 136 - now Bahamas
$ws.Range("A136").Value = "Bahamas"
$ws.Range("B136").Value = 2585
$ws.Range("C136").Value = 39
$ws.Range("D136").Value = 976
$ws.Range("E136").Value = 1550
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 59

# Row 137 - now Estonia
$ws.Range("A137").Value = "Estonia"
$ws.Range("B137").Value = 2564
$ws.Range("C137").Value = 32
$ws.Range("D137").Value = 2195
$ws.Range("E137").Value = 305
$ws.Range("H137").Value = 64

# Row 141 - Trinidad yTobago
$ws.Range("B141").Value = 2347
$ws.Range("C141").Value = 70
$ws.Range("D141").Value = 743
$ws.Range("E141").Value = 1566
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = 38
